# Applies the "New crime data collected" weekly update to the 009 Pct CompStat sheet.
# Header (volume/number + date range) text runs are edited in place via Characters()
# so the shared-string table keeps its original entry/index (matches upstream diff).
# Numeric table cells are written directly; the handful of cells that flip between
# "no data" placeholder text ("0" / "***.*") and real numbers are handled by copying
# format+value from an existing same-styled placeholder/number cell elsewhere on the
# sheet, which lets Excel reuse existing style / shared-string entries instead of
# minting new ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (edit the specific run of shared text, in place) ---
$ws.Range("A8").Characters(21, 2).Text = "16"
$ws.Range("C9").Characters(27, 9).Text = "4/17/2023"
$ws.Range("C9").Characters(47, 9).Text = "4/23/2023"

# --- Cells that were "0"/"***.*" placeholders and now carry a real number ---
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = '#,##0'
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = '#,##0'
$ws.Range("H30").Value = -100
$ws.Range("H30").NumberFormat = '#,##0.0;"-"#,##0.0'

# --- Cells that now have no data and must show the "0"/"***.*" placeholder text ---
# (copy format+value from a same-styled placeholder elsewhere on the sheet so the
# shared string table / style table are reused rather than duplicated)
$ws.Range("C14").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F22").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E23").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("E15").Value = -50
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = -12.5
$ws.Range("M15").Value = 133.333333333333
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -62.5
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -13.333333333333
$ws.Range("I16").Value = 50
$ws.Range("J16").Value = 86
$ws.Range("K16").Value = -41.860465116279
$ws.Range("L16").Value = 28.205128205128
$ws.Range("M16").Value = 13.636363636363
$ws.Range("N16").Value = -81.412639405204
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 50
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 15
$ws.Range("I17").Value = 70
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = 16.666666666666
$ws.Range("L17").Value = 55.555555555555
$ws.Range("M17").Value = 55.555555555555
$ws.Range("N17").Value = -54.545454545454
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 60
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 14.285714285714
$ws.Range("I18").Value = 88
$ws.Range("J18").Value = 106
$ws.Range("K18").Value = -16.981132075471
$ws.Range("L18").Value = -4.347826086956
$ws.Range("M18").Value = 10
$ws.Range("N18").Value = -65.354330708661
$ws.Range("C19").Value = 30
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 150
$ws.Range("F19").Value = 89
$ws.Range("G19").Value = 66
$ws.Range("H19").Value = 34.848484848484
$ws.Range("I19").Value = 307
$ws.Range("J19").Value = 286
$ws.Range("K19").Value = 7.342657342657
$ws.Range("L19").Value = 88.343558282208
$ws.Range("M19").Value = 25.819672131147
$ws.Range("N19").Value = -29.74828375286
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 6
$ws.Range("H20").Value = 200
$ws.Range("I20").Value = 13
$ws.Range("J20").Value = 13
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = -23.529411764705
$ws.Range("M20").Value = -13.333333333333
$ws.Range("N20").Value = -91.612903225806
$ws.Range("C21").Value = 50
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 56.25
$ws.Range("F21").Value = 157
$ws.Range("G21").Value = 128
$ws.Range("H21").Value = 22.65625
$ws.Range("I21").Value = 535
$ws.Range("J21").Value = 561
$ws.Range("K21").Value = -4.634581105169
$ws.Range("L21").Value = 46.174863387978
$ws.Range("M21").Value = 23.842592592592
$ws.Range("N21").Value = -58.430458430458
$ws.Range("C23").Value = 3
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = 30
$ws.Range("I23").Value = 42
$ws.Range("K23").Value = -14.285714285714
$ws.Range("L23").Value = -33.333333333333
$ws.Range("M23").Value = 20
$ws.Range("C24").Value = 37
$ws.Range("D24").Value = 53
$ws.Range("E24").Value = -30.188679245283
$ws.Range("F24").Value = 101
$ws.Range("G24").Value = 193
$ws.Range("H24").Value = -47.668393782383
$ws.Range("I24").Value = 450
$ws.Range("J24").Value = 668
$ws.Range("K24").Value = -32.634730538922
$ws.Range("L24").Value = 79.282868525896
$ws.Range("M24").Value = -4.862579281183
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 39
$ws.Range("H25").Value = -2.564102564102
$ws.Range("I25").Value = 142
$ws.Range("J25").Value = 141
$ws.Range("K25").Value = 0.709219858156
$ws.Range("L25").Value = 65.116279069767
$ws.Range("M25").Value = -0.6993006993
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = -42.857142857142
$ws.Range("I26").Value = 14
$ws.Range("J26").Value = 15
$ws.Range("K26").Value = -6.666666666666
$ws.Range("L26").Value = 16.666666666666
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 50
$ws.Range("J27").Value = 27
$ws.Range("K27").Value = -25.925925925925
$ws.Range("L27").Value = 42.857142857142
$ws.Range("J30").Value = 7

